# Append-edit: refresh fetch timestamp and insert 3 newly discovered
# Lancers listings into the "ランサーズ" worksheet at their correct
# priority-sorted positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-09-23 12:36:21"

# ---------------------------------------------------------------
# 1) Remove the existing hyperlinks up front. Row inserts shift the
#    underlying cell data down but this engine does not relocate the
#    <hyperlink ref="..."> entries that go with them, so we clear them
#    here and recreate every hyperlink (old + new) once the table has
#    its final shape.
# ---------------------------------------------------------------
$ws.Range("F2:F8").Hyperlinks.Delete()

# ---------------------------------------------------------------
# 2) Insert the three new rows at their sorted positions (inserting
#    top-to-bottom so each subsequent index already accounts for the
#    earlier shifts).
# ---------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(8).Insert()
$ws.Rows.Item(10).Insert()

# ---------------------------------------------------------------
# 3) Populate the newly inserted rows.
# ---------------------------------------------------------------

# Row 5: 冠婚葬祭業公式サイトのPHP+MySQLバージョンアップ依頼
$ws.Cells.Item(5, 1).Value = $newTimestamp
$ws.Cells.Item(5, 2).Value = "【急募】冠婚葬祭業公式サイトのPHP+MySQLバージョンアップ依頼"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5385491"
$ws.Cells.Item(5, 7).Value = 100
$ws.Cells.Item(5, 8).Value = "◇MySQL ○PHP"

# Row 8: 仮想通貨トレードの運用とコンサル【1名】のみ募集
$ws.Cells.Item(8, 1).Value = $newTimestamp
$ws.Cells.Item(8, 2).Value = "仮想通貨トレードの運用とコンサル【1名】のみ募集"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5393695"
$ws.Cells.Item(8, 7).Value = 55
$ws.Cells.Item(8, 8).Value = "◆コンサル"

# Row 10: 【相談から実装まで伴走できる方歓迎】介護・福祉×テクノロジー事例収集の仕組みづくり
$ws.Cells.Item(10, 1).Value = $newTimestamp
$ws.Cells.Item(10, 2).Value = "【相談から実装まで伴走できる方歓迎】介護・福祉×テクノロジー事例収集の仕組みづくり"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5398932"
$ws.Cells.Item(10, 7).Value = 18

# ---------------------------------------------------------------
# 4) Refresh the "取得日時" column for every data row (new + old).
# ---------------------------------------------------------------
$ws.Range("A2:A11").Value = $newTimestamp

# ---------------------------------------------------------------
# 5) Recreate the hyperlinks (URL column F) for all 10 data rows and
#    restore the hyperlink style used throughout the sheet.
# ---------------------------------------------------------------
$hyperlinkStyle = $ws.Range("F2").Style

$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5398662")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5398562")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5251319")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5385491")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5398432")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5398772")
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5393695")
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5398657")
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5398932")
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5398736")

$ws.Range("F2:F11").Style = $hyperlinkStyle
